# Update to R4 QA
# - bump IG version (1.0.0 -> 1.1.0), fhirVersion (3.0.1 -> 4.0.0) and the
#   IG canonical URL on the "meta" sheet
# - drop the "STU3/" path segment from every profile StructureDefinition URL
#   on the "profiles" sheet
# - move the active selection/tab around to match the author's last-saved
#   cursor position

$wb = $excel.ActiveWorkbook

# --- meta sheet: version bump -------------------------------------------
$meta = $wb.Worksheets.Item("meta")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B6").Value = "http://hl7.org/fhir/us/davinci-deqm/ImplementationGuide/hl7.fhir.us.davinci-deqm-1.1.0"
$meta.Range("B4").Value = "4.0.0"

# --- profiles sheet: STU3 -> R4 structure definition urls ---------------
$profiles = $wb.Worksheets.Item("profiles")
$profiles.Range("A2").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/devicerequest-deqm"
$profiles.Range("A3").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/medicationadministration-deqm"
$profiles.Range("A4").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/summary-measurereport-deqm"
$profiles.Range("A5").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/datax-measurereport-deqm"
$profiles.Range("A6").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/practitioner-deqm"
$profiles.Range("A7").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/deviceusestatement-deqm"
$profiles.Range("A8").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/organization-deqm"
$profiles.Range("A9").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/coverage-deqm"
$profiles.Range("A10").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/indv-measurereport-deqm"
$profiles.Range("A11").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/medicationrequest-deqm"
$profiles.Range("A12").Value = "!http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/subscription-deqm"

# --- view / selection state ----------------------------------------------
$meta.Activate() | Out-Null
$meta.Range("B5").Select() | Out-Null

$profiles.Activate() | Out-Null
$profiles.Range("A16").Select() | Out-Null
